$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill column M (icon) first for rows 16-18 (type = coins/gems rewards)
$ws.Range("M16").Value = "GR_coins_pack06"
$ws.Range("M17").Value = "GR_gems_pack05"
$ws.Range("M18").Value = "GR_gems_pack05"

# Then fill column N (tid) for rows 16-18
$ws.Range("N16").Value = "TID_SC_NAME_PLURAL"
$ws.Range("N17").Value = "TID_GEM_PLURAL"
$ws.Range("N18").Value = "TID_GEM_PLURAL"

# Update view state: selection moves to L8, and scroll so column A is visible (no topLeftCell override)
$ws.Range("L8").Select()
